$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-08 -> 2023-09-09, i.e. 45177 -> 45178) for every data row
# (rows 2 through 90).
for ($row = 2; $row -le 90; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value = 45178
    }
}
